$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.494.12'
$ws.Range("E2").Value = '  +3.07%  '
$ws.Range("D3").Value = '1.841.79'
$ws.Range("E3").Value = '  +2.43%  '
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.69'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.620'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.43%  '
$ws.Range("E7").Value = '  +0.30%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '43.90'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +11.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.310'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +7.97%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0698'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.77%  '
$ws.Range("E11").Value = '  +2.79%  '
$ws.Range("D12").Value = '2.108.11'
$ws.Range("D13").Value = '1.837.66'
$ws.Range("E13").Value = '  +2.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '11.26'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.66%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.672'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +7.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.71'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +8.14%  '
$ws.Range("D17").Value = '35.493.90'
$ws.Range("E17").Value = '  +3.18%  '
$ws.Range("E18").Value = '  +3.80%  '
$ws.Range("D19").Value = '0.0₃0800'
$ws.Range("E19").Value = '  +4.98%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '244.25'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.32%  '
$ws.Range("E21").Value = '  +8.82%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.63'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +13.98%  '
$ws.Range("E23").Value = '  +0.27%  '
$ws.Range("E24").Value = '  +4.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '171.29'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.45%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.99'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.43%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.80'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.29%  '
$ws.Range("E29").Value = '  +27.75%  '
$ws.Range("E30").Value = '  +0.26%  '
$ws.Range("D31").Value = '3.345.83'
$ws.Range("E31").Value = '  +37.71%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0551'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +7.74%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.09'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.83%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.93'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.40%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.83'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '94.69'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +16.52%  '
$ws.Range("E37").Value = '  +8.69%  '
$ws.Range("E38").Value = '  +7.52%  '
$ws.Range("D39").Value = '1.347.03'
$ws.Range("E39").Value = '  +3.62%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0195'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.83%  '
$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '15.36'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +9.64%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.44'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.17%  '
$ws.Range("E43").Value = '  +7.67%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.26'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.78%  '
$ws.Range("E45").Value = '  +0.80%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.80'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.49%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.26'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +9.78%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0520'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.79%  '
$ws.Range("D49").Value = '2.012.41'
$ws.Range("E49").Value = '  +2.75%  '
$ws.Range("E50").Value = '  +0.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '102.39'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.04%  '
